$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'39.837.50"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = "'2.209.45"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'291.25"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('D6').Value = "'86.79"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.76%  '
$ws.Range('D7').Value = "'0.513"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = "'0.467"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.28%  '
$ws.Range('D10').Value = "'30.23"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.71%  '
$ws.Range('D11').Value = "'0.0778"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.98%  '
$ws.Range('D12').Value = "'50.00"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.51%  '
$ws.Range('E13').Value = '  +2.37%  '
$ws.Range('D14').Value = "'6.42"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('D15').Value = "'2.550.30"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.46%  '
$ws.Range('D16').Value = "'13.72"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.30%  '
$ws.Range('D17').Value = "'2.186.15"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.71%  '
$ws.Range('D18').Value = "'0.729"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.47%  '
$ws.Range('D19').Value = "'39.764.63"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').Value = "'0.0₃0883"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').Value = "'11.13"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.40%  '
$ws.Range('D22').Value = "'5.72"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.48%  '
$ws.Range('D23').Value = "'65.45"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('D24').Value = "'236.90"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').Value = "'2.45"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.61%  '
$ws.Range('E27').Value = '  -0.79%  '
$ws.Range('D28').Value = "'23.11"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.44%  '
$ws.Range('E29').Value = '  +1.08%  '
$ws.Range('D30').Value = "'9.20"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.30%  '
$ws.Range('D31').Value = "'156.80"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.31%  '
$ws.Range('D32').Value = "'31.83"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.21%  '
$ws.Range('D33').Value = "'0.999"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').Value = "'4.94"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').Value = "'0.0708"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.53%  '
$ws.Range('D36').Value = "'2.92"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.19%  '
$ws.Range('E37').Value = '  -1.60%  '
$ws.Range('D38').Value = "'0.111"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.83%  '
$ws.Range('D39').Value = "'0.0979"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.42%  '
$ws.Range('D40').Value = "'1.71"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.46%  '
$ws.Range('D41').Value = "'15.18"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.41%  '
$ws.Range('D42').Value = "'2.111.05"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.20%  '
$ws.Range('D43').Value = "'3.71"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.16%  '
$ws.Range('D44').Value = "'0.0268"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = "'9.97"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.48%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = "'2.09"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.19%  '
$ws.Range('D47').Value = "'17.77"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.59%  '
$ws.Range('D48').Value = "'2.69"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.38%  '
$ws.Range('D49').Value = "'2.422.12"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.43%  '
$ws.Range('E50').Value = '  +2.04%  '
$ws.Range('D51').Value = "'88.43"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.70%  '
